# "Peter-[Update]-create each article analyzer"
# - Update 金融產品 (sheet1) rows 2-3: new product rows (信用卡 / 台新@GoGo卡) with updated counts
# - Update 金融機構 (sheet2) rows 2-11: refreshed institution list + counts
# - Append 8 new topic rows to 主題分類 (sheet3)
# - Add a new worksheet 文章分類 (sheet4) with category / occurrence-count data

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet 1: 金融產品 (Financial products)
# ---------------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("金融產品")

$ws1.Range("A2").Value = "信用卡"
$ws1.Range("B2").Value = 3
$ws1.Range("A3").Value = "台新@GoGo卡"
$ws1.Range("B3").Value = 2
$ws1.Range("A4").Value = "信貸專案"
$ws1.Range("B4").Value = 1
$ws1.Range("A5").Value = "袋鼠金融限定專案"
$ws1.Range("B5").Value = 1
$ws1.Range("A6").Value = "Online貸"
$ws1.Range("B6").Value = 1
$ws1.Range("A7").Value = "信貸卡友專屬方案"
$ws1.Range("B7").Value = 1
$ws1.Range("A8").Value = "數時貸"
$ws1.Range("B8").Value = 1
$ws1.Range("A9").Value = "折扣碼"
$ws1.Range("B9").Value = 1
$ws1.Range("A10").Value = "電子禮券"
$ws1.Range("B10").Value = 1
$ws1.Range("A11").Value = "無線吸塵器"
$ws1.Range("B11").Value = 1

# ---------------------------------------------------------------------------
# Sheet 2: 金融機構 (Financial institutions)
# ---------------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("金融機構")

$ws2.Range("A2").Value = "永豐銀行"
$ws2.Range("B2").Value = 6
$ws2.Range("A3").Value = "玉山銀行"
$ws2.Range("B3").Value = 3
$ws2.Range("A4").Value = "袋鼠金融"
$ws2.Range("B4").Value = 3
$ws2.Range("A5").Value = "台新銀行"
$ws2.Range("B5").Value = 3
$ws2.Range("A6").Value = "富邦銀行"
$ws2.Range("B6").Value = 3
$ws2.Range("A7").Value = "元富證券"
$ws2.Range("B7").Value = 3
$ws2.Range("A8").Value = "渣打銀行"
$ws2.Range("B8").Value = 2
$ws2.Range("A9").Value = "滙豐銀行"
$ws2.Range("B9").Value = 2
$ws2.Range("A10").Value = "凱基銀行"
$ws2.Range("B10").Value = 2
$ws2.Range("A11").Value = "北富銀行"
$ws2.Range("B11").Value = 2

# ---------------------------------------------------------------------------
# Sheet 3: 主題分類 (Topic classification) - append 8 new topics
# ---------------------------------------------------------------------------
$ws3 = $wb.Worksheets.Item("主題分類")

$ws3.Range("A4").Value = "長榮海運的配息政策與股利發放資訊"
$ws3.Range("A5").Value = "耐吉公司及其股價表現與投資建議"
$ws3.Range("A6").Value = "王道信貸與國泰信貸的比較及其他信貸方案介紹"
$ws3.Range("A7").Value = "新戶開戶優惠及投資策略"
$ws3.Range("A8").Value = "生命靈數的性格特質分析"
$ws3.Range("A9").Value = "2025年信用卡推薦及權益變動"
$ws3.Range("A10").Value = "2025年繳稅行事曆及信用卡推薦"
$ws3.Range("A11").Value = "機場旅平險的保障內容及投保管道比較"

# ---------------------------------------------------------------------------
# Sheet 4 (new): 文章分類 (Article classification)
# ---------------------------------------------------------------------------
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$ws4 = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $lastSheet)
$ws4.Name = "文章分類"

# Copy the bold/bordered header formatting used by the other sheets, then
# overwrite the copied values with this sheet's own header text.
$ws1.Range("A1:B1").Copy($ws4.Range("A1:B1"))
$ws4.Range("A1").Value = "分類"
$ws4.Range("B1").Value = "出現次數"

$ws4.Range("A2").Value = "ROO 投資"
$ws4.Range("B2").Value = 1
$ws4.Range("A3").Value = "證券"
$ws4.Range("B3").Value = 1
$ws4.Range("A4").Value = "ROO 貸款"
$ws4.Range("B4").Value = 1
$ws4.Range("A5").Value = "信貸推薦"
$ws4.Range("B5").Value = 1
$ws4.Range("A6").Value = "ROO 時事快訊"
$ws4.Range("B6").Value = 1
$ws4.Range("A7").Value = "熱門話題"
$ws4.Range("B7").Value = 1
$ws4.Range("A8").Value = "ROO 信用卡"
$ws4.Range("B8").Value = 1
$ws4.Range("A9").Value = "信用卡推薦"
$ws4.Range("B9").Value = 1

Write-Output "edit applied"
